$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Split the "E02" error code row into "E02-A" and "E02-B" rows ---
# Row 3 currently holds E02 / "args[1] i.e. sy-uname not specified."
$e02Message = $ws.Range("B3").Value2

# Insert a new blank row right after row 3 (pushes everything below down by one)
$ws.Rows.Item(4).Insert()

$ws.Range("A3").Value = "E02-A"
$ws.Range("B3").Value = $e02Message

$ws.Range("A4").Value = "E02-B"
$ws.Range("B4").Value = $e02Message

# --- Split the "E12" error code row into "E12-A" and "E12-B" rows ---
# After the insert above, the old E12 row (previously row 8) is now row 9
$e12Message = $ws.Range("B9").Value2

# Insert a new blank row right after row 9
$ws.Rows.Item(10).Insert()

$ws.Range("A9").Value = "E12-A"
$ws.Range("B9").Value = $e12Message

$ws.Range("A10").Value = "E12-B"
$ws.Range("B10").Value = $e12Message

# --- Fix up the view so it no longer shows the old scrolled/selected state ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A11").Select()
